$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")
$ws.Rows.Item(35).Insert()
$ws.Range("R35").Value = "beneficiary saravanan"
$ws.Range("S35").Value = "2024-09-09 14:52:20"
